$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 40: Sani-T-10 Sanitizer
Set-TextValue $ws.Cells.Item(40, 1) ""
Set-TextValue $ws.Cells.Item(40, 2) "Sani-T-10 Sanitizer"
Set-TextValue $ws.Cells.Item(40, 3) "1"
Set-TextValue $ws.Cells.Item(40, 4) "107.36"
Set-TextValue $ws.Cells.Item(40, 5) "107.36"

# Row 41: Tamper Evident - 12oz Bowl (Smoothie)
Set-TextValue $ws.Cells.Item(41, 1) ""
Set-TextValue $ws.Cells.Item(41, 2) "Tamper Evident - 12oz Bowl (Smoothie)"
Set-TextValue $ws.Cells.Item(41, 3) "1"
Set-TextValue $ws.Cells.Item(41, 4) "0.00"
Set-TextValue $ws.Cells.Item(41, 5) "0.00"
